$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.719.27"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.601.58"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'211.57"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'19.70"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.826.93"
$ws.Range("D13").Value = "1.605.64"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "26.694.06"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "'210.24"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "'8.96"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "'144.18"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'15.38"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'0.0510"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "1.295.34"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("E38").Value = "  +7.07%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "1.739.08"
$ws.Range("D46").Value = "'90.81"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "'0.102"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = "  +0.07%  "
